# Updated symbol list on Wed Dec 21 05:44:02 UTC 2022 with GitHub Actions
#
# Refresh the "Price" (D) and some "Volume(1h)" (E) columns with the latest
# scraped values, and re-sort the three rows around BKEXToken / CEJI /
# KickToken to their new ranking order.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $value) {
    # The Price column stores numeric-looking quotes as literal text
    # (t="inlineStr" in the original file). Force the cell to a text
    # format before assigning so Excel doesn't silently reinterpret
    # "248.72" as the number 248.72.
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $value
}

# --- simple price refreshes -------------------------------------------------
Set-TextValue "D2"  "248.72"
Set-TextValue "D3"  "22.43"
Set-TextValue "D4"  "5.332"
Set-TextValue "D5"  "0.05687"
Set-TextValue "D6"  "3.401"
Set-TextValue "D7"  "6.332"
Set-TextValue "D8"  "0.8130"
Set-TextValue "D9"  "0.9168"
Set-TextValue "D10" "0.1410"
Set-TextValue "D11" "0.07438"
Set-TextValue "D12" "0.03101"
Set-TextValue "D13" "0.03018"
Set-TextValue "D14" "0.09359"
Set-TextValue "D15" "3.718"
Set-TextValue "D16" "0.001586"
Set-TextValue "D18" "0.01827"

Set-TextValue "D19" "0.0005789"
$ws.Range("E19").Value = "18OneONEWorstin24h"

Set-TextValue "D20" "0.006443"
Set-TextValue "D22" "0.001024"
Set-TextValue "D24" "3.697"
Set-TextValue "D25" "2.136"
Set-TextValue "D27" "0.1298"
Set-TextValue "D40" "0.03978"

# --- rows 41-43 reshuffled: KickToken, BKEXToken, CEJI ----------------------
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
Set-TextValue "D41" "0.006889"
$ws.Range("E41").Value = "40KickTokenKICK"

$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue "D42" "0.1065"
$ws.Range("E42").Value = "41BKEXTokenBKK"

$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue "D43" "0.002709"
$ws.Range("E43").Value = "42CEJICEJI"

# --- remaining price refreshes ----------------------------------------------
Set-TextValue "D44" "0.007445"
Set-TextValue "D45" "0.00005889"
Set-TextValue "D47" "0.4999"
Set-TextValue "D48" "0.2155"
